# Move the "Straight Connector 84" connector shape (id=85) earlier in the
# slide's z-order so that it sits directly after "Straight Arrow Connector 43"
# (id=44) / directly before the "Rectangle 62" shape (id=53), instead of its
# previous spot directly after "TextBox 73" (id=74).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connector = $null
$anchor = $null
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Straight Connector 84") { $connector = $shp }
    if ($shp.Name -eq "Straight Arrow Connector 43") { $anchor = $shp }
}

if ($connector -eq $null) {
    throw "Could not find shape 'Straight Connector 84'"
}
if ($anchor -eq $null) {
    throw "Could not find shape 'Straight Arrow Connector 43'"
}

# Send the connector backward (toward the back of the z-order) until it is
# immediately in front of the anchor shape, i.e. one position after it.
$steps = $connector.ZOrderPosition - $anchor.ZOrderPosition - 1
for ($i = 0; $i -lt $steps; $i++) {
    $connector.ZOrder(3)   # msoSendBackward
}
